$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 11.7268977612736
$ws.Cells.Item(2, 3).Value = 9.79233391983235
$ws.Cells.Item(2, 5).Value = 16.6512875608202
$ws.Cells.Item(2, 6).Value = 16.86991607391245
$ws.Cells.Item(2, 7).Value = 15.6196625311847
$ws.Cells.Item(2, 8).Value = 10.69305984598574
$ws.Cells.Item(2, 9).Value = 14.52566069640675
$ws.Cells.Item(2, 14).Value = 15.23256114932579
$ws.Cells.Item(2, 15).Value = 14.51287010159328
$ws.Cells.Item(3, 2).Value = 11.03122198599277
$ws.Cells.Item(3, 3).Value = 9.393387448607269
$ws.Cells.Item(3, 5).Value = 15.6978117464415
$ws.Cells.Item(3, 6).Value = 15.89584955866815
$ws.Cells.Item(3, 7).Value = 15.47850989842839
$ws.Cells.Item(3, 8).Value = 10.73371870449598
$ws.Cells.Item(3, 9).Value = 14.63123272557883
$ws.Cells.Item(3, 14).Value = 15.24461718816648
$ws.Cells.Item(3, 15).Value = 14.55279137027828
$ws.Cells.Item(4, 2).Value = 10.57892014016331
$ws.Cells.Item(4, 3).Value = 9.138082252656059
$ws.Cells.Item(4, 5).Value = 15.08648893994878
$ws.Cells.Item(4, 6).Value = 15.26997757108489
$ws.Cells.Item(4, 7).Value = 15.40075399083871
$ws.Cells.Item(4, 8).Value = 10.76102292764742
$ws.Cells.Item(4, 9).Value = 14.69997423949106
$ws.Cells.Item(4, 14).Value = 15.25369081518996
$ws.Cells.Item(4, 15).Value = 14.58222275350954
$ws.Cells.Item(5, 2).Value = 10.38832467662627
$ws.Cells.Item(5, 3).Value = 9.031540154861695
$ws.Cells.Item(5, 5).Value = 14.83112576560417
$ws.Cells.Item(5, 6).Value = 15.008197319934
$ws.Cells.Item(5, 7).Value = 15.3713462812162
$ws.Cells.Item(5, 8).Value = 10.77273643593765
$ws.Cells.Item(5, 9).Value = 14.72897144460455
$ws.Cells.Item(5, 14).Value = 15.25780975829413
$ws.Cells.Item(5, 15).Value = 14.59544637573338
$ws.Cells.Item(6, 2).Value = 10.35629884973616
$ws.Cells.Item(6, 3).Value = 9.013700914603874
$ws.Cells.Item(6, 5).Value = 14.78835455969835
$ws.Cells.Item(6, 6).Value = 14.96433081551589
$ws.Cells.Item(6, 7).Value = 15.36660172255279
$ws.Cells.Item(6, 8).Value = 10.77471685388701
$ws.Cells.Item(6, 9).Value = 14.73384583537939
$ws.Cells.Item(6, 14).Value = 15.2585191876214
$ws.Cells.Item(6, 15).Value = 14.59771623001476
$ws.Cells.Item(7, 2).Value = 10.57637505445486
$ws.Cells.Item(7, 3).Value = 9.136655382276267
$ws.Cells.Item(7, 5).Value = 15.08306991642536
$ws.Cells.Item(7, 6).Value = 15.26647399323133
$ws.Cells.Item(7, 7).Value = 15.40034811862107
$ws.Cells.Item(7, 8).Value = 10.76117852610564
$ws.Cells.Item(7, 9).Value = 14.70036132136313
$ws.Cells.Item(7, 14).Value = 15.25374465713749
$ws.Cells.Item(7, 15).Value = 14.58239612089481
$ws.Cells.Item(8, 2).Value = 11.49227329408074
$ws.Cells.Item(8, 3).Value = 9.656980998844524
$ws.Cells.Item(8, 5).Value = 16.32804466967644
$ws.Cells.Item(8, 6).Value = 16.5399640634477
$ws.Cells.Item(8, 7).Value = 15.56916786608193
$ws.Cells.Item(8, 8).Value = 10.70659257075207
$ws.Cells.Item(8, 9).Value = 14.56124762220964
$ws.Cells.Item(8, 14).Value = 15.23637191188535
$ws.Cells.Item(8, 15).Value = 14.52560984909625
$ws.Cells.Item(9, 2).Value = 13.0873429154736
$ws.Cells.Item(9, 3).Value = 10.59149559864832
$ws.Cells.Item(9, 5).Value = 18.68939558627889
$ws.Cells.Item(9, 6).Value = 19.00274580682531
$ws.Cells.Item(9, 7).Value = 15.96896040893071
$ws.Cells.Item(9, 8).Value = 10.61817822578057
$ws.Cells.Item(9, 9).Value = 14.31960410739308
$ws.Cells.Item(9, 14).Value = 15.21551368756914
$ws.Cells.Item(9, 15).Value = 14.45357658611584
$ws.Cells.Item(10, 2).Value = 14.13549404846668
$ws.Cells.Item(10, 3).Value = 11.22159403717808
$ws.Cells.Item(10, 5).Value = 20.34839710100202
$ws.Cells.Item(10, 6).Value = 20.67494806633232
$ws.Cells.Item(10, 7).Value = 16.30151156664
$ws.Cells.Item(10, 8).Value = 10.56466557514306
$ws.Cells.Item(10, 9).Value = 14.16114311944878
$ws.Cells.Item(10, 14).Value = 15.20817393855715
$ws.Cells.Item(10, 15).Value = 14.42499002443218
$ws.Cells.Item(11, 2).Value = 14.58534766749324
$ws.Cells.Item(11, 3).Value = 11.49530076729246
$ws.Cells.Item(11, 5).Value = 21.06034307139853
$ws.Cells.Item(11, 6).Value = 21.3917225636224
$ws.Cells.Item(11, 7).Value = 16.46046742105951
$ws.Cells.Item(11, 8).Value = 10.54282529428894
$ws.Cells.Item(11, 9).Value = 14.09321657196579
$ws.Cells.Item(11, 14).Value = 15.20655338590133
$ws.Cells.Item(11, 15).Value = 14.41733557331771
$ws.Cells.Item(12, 2).Value = 14.75181633261379
$ws.Cells.Item(12, 3).Value = 11.59704214871203
$ws.Cells.Item(12, 5).Value = 21.32382868579417
$ws.Cells.Item(12, 6).Value = 21.65686569030329
$ws.Cells.Item(12, 7).Value = 16.52169476701663
$ws.Cells.Item(12, 8).Value = 10.53491639408536
$ws.Cells.Item(12, 9).Value = 14.06809448877751
$ws.Cells.Item(12, 14).Value = 15.20618542533379
$ws.Cells.Item(12, 15).Value = 14.41521072529082
$ws.Cells.Item(13, 2).Value = 14.71613697704697
$ws.Cells.Item(13, 3).Value = 11.57521572613386
$ws.Cells.Item(13, 5).Value = 21.26735368677254
$ws.Cells.Item(13, 6).Value = 21.60004134736742
$ws.Cells.Item(13, 7).Value = 16.50846350361657
$ws.Cells.Item(13, 8).Value = 10.53660361201894
$ws.Cells.Item(13, 9).Value = 14.07347824567568
$ws.Cells.Item(13, 14).Value = 15.20625376694213
$ws.Cells.Item(13, 15).Value = 14.41563387392817
$ws.Cells.Item(14, 2).Value = 14.59912094228819
$ws.Cells.Item(14, 3).Value = 11.50370951879382
$ws.Cells.Item(14, 5).Value = 21.08214251211878
$ws.Cells.Item(14, 6).Value = 21.4136618050453
$ws.Cells.Item(14, 7).Value = 16.46548431323675
$ws.Cells.Item(14, 8).Value = 10.54216736927917
$ws.Cells.Item(14, 9).Value = 14.09113772001547
$ws.Cells.Item(14, 14).Value = 15.20651819734519
$ws.Cells.Item(14, 15).Value = 14.417145231064
$ws.Cells.Item(15, 2).Value = 14.52693972733271
$ws.Cells.Item(15, 3).Value = 11.45966051070594
$ws.Cells.Item(15, 5).Value = 20.96790023029221
$ws.Cells.Item(15, 6).Value = 21.29868154950795
$ws.Cells.Item(15, 7).Value = 16.43929088836378
$ws.Cells.Item(15, 8).Value = 10.54562246230083
$ws.Cells.Item(15, 9).Value = 14.10203289550018
$ws.Cells.Item(15, 14).Value = 15.20671212494251
$ws.Cells.Item(15, 15).Value = 14.41817186356781
$ws.Cells.Item(16, 2).Value = 14.10555212023672
$ws.Cells.Item(16, 3).Value = 11.20344237627911
$ws.Cells.Item(16, 5).Value = 20.3010125364109
$ws.Cells.Item(16, 6).Value = 20.62722412089977
$ws.Cells.Item(16, 7).Value = 16.29127201063987
$ws.Cells.Item(16, 8).Value = 10.56614339314171
$ws.Cells.Item(16, 9).Value = 14.16566614411996
$ws.Cells.Item(16, 14).Value = 15.2083143053872
$ws.Cells.Item(16, 15).Value = 14.42559830568041
$ws.Cells.Item(17, 2).Value = 13.84013559190855
$ws.Cells.Item(17, 3).Value = 11.04291423207018
$ws.Cells.Item(17, 5).Value = 19.88097508528089
$ws.Cells.Item(17, 6).Value = 20.20408069597325
$ws.Cells.Item(17, 7).Value = 16.20238476296524
$ws.Cells.Item(17, 8).Value = 10.5793745362031
$ws.Cells.Item(17, 9).Value = 14.20576966694501
$ws.Cells.Item(17, 14).Value = 15.20973632749838
$ws.Cells.Item(17, 15).Value = 14.43152774809081
$ws.Cells.Item(18, 2).Value = 13.68493782922282
$ws.Cells.Item(18, 3).Value = 10.94936808520677
$ws.Cells.Item(18, 5).Value = 19.63535557912909
$ws.Cells.Item(18, 6).Value = 19.95656407809801
$ws.Cells.Item(18, 7).Value = 16.1519861159676
$ws.Cells.Item(18, 8).Value = 10.58722022795361
$ws.Cells.Item(18, 9).Value = 14.2292273915987
$ws.Cells.Item(18, 14).Value = 15.2107160610428
$ws.Cells.Item(18, 15).Value = 14.43544152239063
$ws.Cells.Item(19, 2).Value = 13.63195498919732
$ws.Cells.Item(19, 3).Value = 10.91748789379916
$ws.Cells.Item(19, 5).Value = 19.55150070289924
$ws.Cells.Item(19, 6).Value = 19.87204792380568
$ws.Cells.Item(19, 7).Value = 16.13504892936204
$ws.Cells.Item(19, 8).Value = 10.58991703972937
$ws.Cells.Item(19, 9).Value = 14.23723690240814
$ws.Cells.Item(19, 14).Value = 15.21107561686636
$ws.Cells.Item(19, 15).Value = 14.4368529580251
$ws.Cells.Item(20, 2).Value = 13.86865227378903
$ws.Cells.Item(20, 3).Value = 11.06012879343201
$ws.Cells.Item(20, 5).Value = 19.92610512641154
$ws.Cells.Item(20, 6).Value = 20.24955283636154
$ws.Cells.Item(20, 7).Value = 16.21177222168914
$ws.Cells.Item(20, 8).Value = 10.57794167439832
$ws.Cells.Item(20, 9).Value = 14.20146007054665
$ws.Cells.Item(20, 14).Value = 15.20956821251198
$ws.Cells.Item(20, 15).Value = 14.43084442320589
$ws.Cells.Item(21, 2).Value = 14.63359674022021
$ws.Cells.Item(21, 3).Value = 11.52476467072559
$ws.Cells.Item(21, 5).Value = 21.13670916035013
$ws.Cells.Item(21, 6).Value = 21.46857628470577
$ws.Cells.Item(21, 7).Value = 16.4780808437746
$ws.Cells.Item(21, 8).Value = 10.54052333396801
$ws.Cells.Item(21, 9).Value = 14.08593439540717
$ws.Cells.Item(21, 14).Value = 15.20643387059965
$ws.Cells.Item(21, 15).Value = 14.41668027799007
$ws.Cells.Item(22, 2).Value = 15.11091037792833
$ws.Cells.Item(22, 3).Value = 11.8173083426407
$ws.Cells.Item(22, 5).Value = 21.8922937051427
$ws.Cells.Item(22, 6).Value = 22.22866616901552
$ws.Cells.Item(22, 7).Value = 16.65812305655683
$ws.Cells.Item(22, 8).Value = 10.51817635415173
$ws.Cells.Item(22, 9).Value = 14.01393143078454
$ws.Cells.Item(22, 14).Value = 15.2058168890852
$ws.Cells.Item(22, 15).Value = 14.41193464947234
$ws.Cells.Item(23, 2).Value = 14.85822960924761
$ws.Cells.Item(23, 3).Value = 11.66220346693967
$ws.Cells.Item(23, 5).Value = 21.49227183169648
$ws.Cells.Item(23, 6).Value = 21.82633154458858
$ws.Cells.Item(23, 7).Value = 16.56150635647986
$ws.Cells.Item(23, 8).Value = 10.5299099387007
$ws.Cells.Item(23, 9).Value = 14.05203974770672
$ws.Cells.Item(23, 14).Value = 15.20601567355275
$ws.Cells.Item(23, 15).Value = 14.41405338611882
$ws.Cells.Item(24, 2).Value = 13.85576799948579
$ws.Cells.Item(24, 3).Value = 11.05235000006941
$ws.Cells.Item(24, 5).Value = 19.90571471559496
$ws.Cells.Item(24, 6).Value = 20.22900810905287
$ws.Cells.Item(24, 7).Value = 16.2075259539403
$ws.Cells.Item(24, 8).Value = 10.57858872704876
$ws.Cells.Item(24, 9).Value = 14.20340718800242
$ws.Cells.Item(24, 14).Value = 15.20964371195254
$ws.Cells.Item(24, 15).Value = 14.43115178206839
$ws.Cells.Item(25, 2).Value = 12.67750883299381
$ws.Cells.Item(25, 3).Value = 10.34834513028258
$ws.Cells.Item(25, 5).Value = 18.04032842806736
$ws.Cells.Item(25, 6).Value = 18.34778573295695
$ws.Cells.Item(25, 7).Value = 15.85374410594886
$ws.Cells.Item(25, 8).Value = 10.64009308548783
$ws.Cells.Item(25, 9).Value = 14.38163164102951
$ws.Cells.Item(25, 14).Value = 15.21974905795872
$ws.Cells.Item(25, 15).Value = 14.46881386243014
